$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "34.339.31"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +12.15%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.825.27"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +9.19%  "

$ws.Range("E4").Value = "  +0.01%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "229.93"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +5.05%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.578"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +9.53%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -0.03%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "31.51"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +8.13%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "46.59"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +5.53%  "

$ws.Range("E10").Value = "  +9.92%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0680"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +6.42%  "

$ws.Range("E12").Value = "  +3.30%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "2.087.15"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +9.21%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "1.816.45"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +8.49%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.650"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +8.15%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "34.351.02"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +12.23%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "10.32"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +2.77%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "4.30"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +7.37%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "70.65"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +7.26%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "258.56"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +6.83%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0757"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +5.46%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +0.00%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "10.69"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +7.32%  "

$ws.Range("E24").Value = "  +2.99%  "

$ws.Range("E25").Value = "  +3.46%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "159.30"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.12%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "16.80"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +6.50%  "

$ws.Range("E28").Value = "  +5.81%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.18"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +7.65%  "

$ws.Range("E30").Value = "  -0.03%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "3.91"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +13.34%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.0526"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +6.89%  "

$ws.Range("E33").Value = "  +6.89%  "

$ws.Range("E34").Value = "  +8.10%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.532.56"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +2.29%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.81"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +2.59%  "

$ws.Range("E37").Value = "  +6.70%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.639"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +7.17%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.0190"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +7.20%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "84.45"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +1.05%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "2.81"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +5.19%  "

$ws.Range("E42").Value = "  +3.06%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.914"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +9.35%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.12"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +5.57%  "

$ws.Range("E45").Value = "  +5.71%  "

$ws.Range("E46").Value = "  +6.27%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.983.46"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +9.83%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "5.86"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +5.88%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "12.26"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +20.10%  "

$ws.Range("E50").Value = "  -0.01%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "51.87"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +4.05%  "
